$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 16, 3, 4),
    @(3, 16, 4, 4),
    @(7, 14, 5, 6),
    @(4, 14, 3, 6),
    @(4, 16, 5, 4),
    @(1, 5, 5, 15),
    @(4, 12, 5, 8),
    @(5, 12, 2, 8),
    @(3, 8, 4, 12),
    @(4, 12, 2, 8),
    @(4, 8, 3, 12),
    @(5, 8, 7, 12),
    @(4, 8, 3, 12),
    @(6, 4, 5, 16),
    @(4, 6, 3, 14)
)

$startRow = 951
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

$ws.Range("A966").Select()
